$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update raw input values (row 6 and row 8) — dependent formulas
# (H6, I6, K6, L6, H8, K8, etc.) recalculate automatically.
$ws.Range("C6").Value = 1350
$ws.Range("D6").Value = 1550
$ws.Range("D8").Value = 1250

# Move the active selection from F12 to E17.
$ws.Range("E17").Select()
